# Update data: 11 March 2022
# Adds the new monthly observation (date serial 44593 = 1-Feb-2022) to both
# the "Canada" sheet and the "Province" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Canada" -> new row 27
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Canada")

$ws1.Range("A27").Value = 44593
$ws1.Range("A27").NumberFormat = "d-mmm-yy"
$ws1.Range("B27").Value = "Canada"
$ws1.Range("B27").NumberFormat = "d-mmm-yy"
$ws1.Range("D27").Value = 1135.5
$ws1.Range("E27").Value = 1177.2
$ws1.Range("C27").Formula = "=(D27-E27)/E27*100"

# ---------------------------------------------------------------------
# Sheet "Province" -> new rows 252-261 (one per province, same date)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Province")

$provinces = @(
    @{row=252; name="Newfoundland & Labrador"; d=32.200000000000003; e=32.200000000000003; first=$true},
    @{row=253; name="Prince Edward Island";     d=8.3000000000000007; e=8.5;                first=$false},
    @{row=254; name="Nova Scotia";              d=33.6;               e=32.9;               first=$false},
    @{row=255; name="New Brunswick";            d=31.1;               e=34.4;               first=$false},
    @{row=256; name="Quebec";                   d=207.2;              e=246;                first=$false},
    @{row=257; name="Ontario";                  d=451.4;              e=445.8;              first=$false},
    @{row=258; name="Manitoba";                 d=34.1;               e=37.200000000000003; first=$false},
    @{row=259; name="Saskatchewan";              d=28.7;               e=35.1;               first=$false},
    @{row=260; name="Alberta";                  d=169;                e=178.8;              first=$false},
    @{row=261; name="British Columbia";         d=140;                e=126.2;              first=$false}
)

foreach ($p in $provinces) {
    $r = $p.row
    $ws2.Range("A$r").Value = 44593
    $ws2.Range("A$r").NumberFormat = "d-mmm-yy"
    $ws2.Range("B$r").Value = $p.name
    if ($p.first) {
        $ws2.Range("B$r").NumberFormat = "d-mmm-yy"
    }
    $ws2.Range("D$r").Value = $p.d
    $ws2.Range("E$r").Value = $p.e
    $ws2.Range("C$r").Formula = "=(D$r-E$r)/E$r*100"
}

# ---------------------------------------------------------------------
# Restore view/selection state to match the new data extent.
# (Province stays the active/selected sheet, as in the source file.)
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A27").Select() | Out-Null

$ws2.Activate()
$ws2.Range("D262").Select() | Out-Null
